# Auto-generated edit script: refresh market-price derived columns (H-N)
# across multiple crafting-class sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 223.6923
$ws.Range("I9").Value = 295.33334
$ws.Range("K9").Value = 295.33334
$ws.Range("M9").Value = -126.33334
$ws.Range("H19").Value = 488.33334
$ws.Range("I19").Value = 565.6
$ws.Range("J19").Value = 391.75
$ws.Range("K19").Value = 565.6
$ws.Range("L19").Value = 391.75
$ws.Range("M19").Value = -390.6
$ws.Range("N19").Value = -741.75
$ws.Range("H28").Value = 1023.9583
$ws.Range("I28").Value = 382.5
$ws.Range("K28").Value = 382.5
$ws.Range("M28").Value = 102.5
$ws.Range("H43").Value = 10630.7
$ws.Range("I43").Value = 5355.5
$ws.Range("J43").Value = 14147.5
$ws.Range("K43").Value = 5355.5
$ws.Range("L43").Value = 14147.5
$ws.Range("M43").Value = -5286.5
$ws.Range("N43").Value = -14285.5
$ws.Range("H62").Value = 4638.353
$ws.Range("I62").Value = 4353.273
$ws.Range("J62").Value = 5161
$ws.Range("K62").Value = 4353.273
$ws.Range("L62").Value = 5161
$ws.Range("M62").Value = -3729.273
$ws.Range("N62").Value = -6409
$ws.Range("H65").Value = 4638.353
$ws.Range("I65").Value = 4353.273
$ws.Range("J65").Value = 5161
$ws.Range("K65").Value = 21766.365
$ws.Range("L65").Value = 25805
$ws.Range("M65").Value = -18646.365
$ws.Range("N65").Value = -32045
$ws.Range("H116").Value = 4800
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4800
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4800
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11684
$ws.Range("H125").Value = 2215.4
$ws.Range("I125").Value = 1358
$ws.Range("J125").Value = 2527.182
$ws.Range("K125").Value = 12222
$ws.Range("L125").Value = 22744.638
$ws.Range("M125").Value = -9762
$ws.Range("N125").Value = -27664.638
$ws.Range("H131").Value = 878.5714
$ws.Range("I131").Value = 878.5714
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2635.7142
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2404.2858
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2602.739
$ws.Range("I132").Value = 2266.5
$ws.Range("K132").Value = 6799.5
$ws.Range("M132").Value = -4269.5
$ws.Range("H135").Value = 63396.625
$ws.Range("I135").Value = 840
$ws.Range("K135").Value = 7560
$ws.Range("M135").Value = -5025
$ws.Range("H137").Value = 2242.7727
$ws.Range("I137").Value = 2056.2666
$ws.Range("J137").Value = 2642.4285
$ws.Range("K137").Value = 6168.7998
$ws.Range("L137").Value = 7927.2855
$ws.Range("M137").Value = -3618.7998
$ws.Range("N137").Value = -13027.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2445.889
$ws.Range("I2").Value = 2002.25
$ws.Range("J2").Value = 5995
$ws.Range("K2").Value = 2002.25
$ws.Range("L2").Value = 5995
$ws.Range("M2").Value = -1889.25
$ws.Range("N2").Value = -6221
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H61").Value = 2349.9678
$ws.Range("I61").Value = 2019.1666
$ws.Range("K61").Value = 2019.1666
$ws.Range("M61").Value = -1807.1666
$ws.Range("H102").Value = 2777.111
$ws.Range("I102").Value = 1535.2174
$ws.Range("K102").Value = 1535.2174
$ws.Range("M102").Value = 86.7826
$ws.Range("H116").Value = 2445.889
$ws.Range("I116").Value = 2002.25
$ws.Range("J116").Value = 5995
$ws.Range("K116").Value = 2002.25
$ws.Range("L116").Value = 5995
$ws.Range("M116").Value = 291.75
$ws.Range("N116").Value = -10583
$ws.Range("H122").Value = 2634
$ws.Range("I122").Value = 2630.1365
$ws.Range("J122").Value = 2648.1667
$ws.Range("K122").Value = 7890.4095
$ws.Range("L122").Value = 7944.500100000001
$ws.Range("M122").Value = -5440.4095
$ws.Range("N122").Value = -12844.5001
$ws.Range("H132").Value = 2345.3235
$ws.Range("I132").Value = 1368.1852
$ws.Range("J132").Value = 6114.2856
$ws.Range("K132").Value = 4104.5556
$ws.Range("L132").Value = 18342.8568
$ws.Range("M132").Value = -1574.5556
$ws.Range("N132").Value = -23402.8568
$ws.Range("H136").Value = 2349.9678
$ws.Range("I136").Value = 2019.1666
$ws.Range("K136").Value = 6057.4998
$ws.Range("M136").Value = -3507.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2445.889
$ws.Range("I3").Value = 2002.25
$ws.Range("J3").Value = 5995
$ws.Range("K3").Value = 2002.25
$ws.Range("L3").Value = 5995
$ws.Range("M3").Value = -1888.25
$ws.Range("N3").Value = -6223
$ws.Range("H26").Value = 23298
$ws.Range("I26").Value = 21492.75
$ws.Range("J26").Value = 30519
$ws.Range("K26").Value = 21492.75
$ws.Range("L26").Value = 30519
$ws.Range("M26").Value = -21200.75
$ws.Range("N26").Value = -31103

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11701.459
$ws.Range("I31").Value = 2582.7334
$ws.Range("K31").Value = 2582.7334
$ws.Range("M31").Value = -2287.7334
$ws.Range("H34").Value = 11701.459
$ws.Range("I34").Value = 2582.7334
$ws.Range("K34").Value = 2582.7334
$ws.Range("M34").Value = -2380.7334
$ws.Range("H58").Value = 3474.2632
$ws.Range("I58").Value = 2779
$ws.Range("J58").Value = 4100
$ws.Range("K58").Value = 2779
$ws.Range("L58").Value = 4100
$ws.Range("M58").Value = -2576
$ws.Range("N58").Value = -4506
$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("H122").Value = 146370
$ws.Range("I122").Value = 146370
$ws.Range("K122").Value = 439110
$ws.Range("M122").Value = -436660
$ws.Range("H134").Value = 2736.647
$ws.Range("I134").Value = 2736.647
$ws.Range("K134").Value = 8209.940999999999
$ws.Range("M134").Value = -5674.940999999999
$ws.Range("H136").Value = 3474.2632
$ws.Range("I136").Value = 2779
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 8337
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -5787
$ws.Range("N136").Value = -17400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 65
$ws.Range("I13").Value = 30
$ws.Range("K13").Value = 90
$ws.Range("M13").Value = 78

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6031.35
$ws.Range("J70").Value = 5936.6665
$ws.Range("L70").Value = 5936.6665
$ws.Range("N70").Value = -6476.6665
$ws.Range("H73").Value = 6031.35
$ws.Range("J73").Value = 5936.6665
$ws.Range("L73").Value = 5936.6665
$ws.Range("N73").Value = -7808.6665
$ws.Range("H80").Value = 2496.5625
$ws.Range("I80").Value = 2017.5652
$ws.Range("J80").Value = 3720.6667
$ws.Range("K80").Value = 2017.5652
$ws.Range("L80").Value = 3720.6667
$ws.Range("M80").Value = -1019.5652
$ws.Range("N80").Value = -5716.6667
$ws.Range("H83").Value = 2496.5625
$ws.Range("I83").Value = 2017.5652
$ws.Range("J83").Value = 3720.6667
$ws.Range("K83").Value = 10087.826
$ws.Range("L83").Value = 18603.3335
$ws.Range("M83").Value = -5095.826000000001
$ws.Range("N83").Value = -28587.3335
$ws.Range("H102").Value = 1910.7916
$ws.Range("I102").Value = 1706.9048
$ws.Range("K102").Value = 1706.9048
$ws.Range("M102").Value = -84.90480000000002
$ws.Range("H132").Value = 4603.3423
$ws.Range("I132").Value = 4040.4644
$ws.Range("J132").Value = 6179.4
$ws.Range("K132").Value = 12121.3932
$ws.Range("L132").Value = 18538.2
$ws.Range("M132").Value = -9591.393199999999
$ws.Range("N132").Value = -23598.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7219.6924
$ws.Range("I40").Value = 3622
$ws.Range("K40").Value = 3622
$ws.Range("M40").Value = -3486
$ws.Range("H61").Value = 45643
$ws.Range("I61").Value = 44542.793
$ws.Range("J61").Value = 49188.11
$ws.Range("K61").Value = 44542.793
$ws.Range("L61").Value = 49188.11
$ws.Range("M61").Value = -44340.793
$ws.Range("N61").Value = -49592.11
$ws.Range("H81").Value = 25181
$ws.Range("J81").Value = 25181
$ws.Range("L81").Value = 25181
$ws.Range("N81").Value = -27177
$ws.Range("H84").Value = 25181
$ws.Range("J84").Value = 25181
$ws.Range("L84").Value = 75543
$ws.Range("N84").Value = -85527
$ws.Range("H113").Value = 45643
$ws.Range("I113").Value = 44542.793
$ws.Range("J113").Value = 49188.11
$ws.Range("K113").Value = 44542.793
$ws.Range("L113").Value = 49188.11
$ws.Range("M113").Value = -42372.793
$ws.Range("N113").Value = -53528.11
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -125060
